$wb = $excel.ActiveWorkbook
$wsWarm = $wb.Worksheets.Item("warm")
$wsNever = $wb.Worksheets.Item("never")

# --- Extend the "theenglish" query table on the "warm" sheet by 6 rows (A479:C484) ---
$lo = $wsWarm.ListObjects.Item(1)
for ($i = 0; $i -lt 6; $i++) {
    $lo.ListRows.Add() | Out-Null
}

# --- Fill in the new rows. Values are written in the same order they were first
#     typed by the original author so new shared strings line up with the source file. ---
$wsNever.Range('A8').Value = 'When did you last win first price in the christmas lottery?'
$wsNever.Range('A9').Value = 'When did you last win an archery contest?'
$wsNever.Range('A10').Value = 'When did you last sit down on a chocolate cake?'
$wsNever.Range('B10').Value = 'I''ve never sat down on a chocolate cake'
$wsNever.Range('B8').Value = 'I''ve never won first price in the christmas lottery'
$wsNever.Range('B9').Value = 'I''ve never won an archery contest'
$wsWarm.Range('B479').Value = 'It catches people off guard'
$wsWarm.Range('A479').Value = 'esto pilla a la gente desprevenida( por sorpresa)'
$wsWarm.Range('C479').Value = 'expressions'
$wsNever.Range('B13').Value = 'I''ve never been in a Malasyan pub'
$wsNever.Range('A15').Value = 'When did you last build an igloo?'
$wsNever.Range('B15').Value = 'I''ve never built an Igloo'
$wsNever.Range('A11').Value = 'When did you last seleep under a bridge?'
$wsNever.Range('B11').Value = 'I''ve never slept under a bridge.'
$wsNever.Range('A12').Value = 'When did you last fly in a hot air balloon?'
$wsNever.Range('B12').Value = 'I''ve never flown in a hot air balloon'
$wsNever.Range('A13').Value = 'When was the last time you were in a Malasyan pub?'
$wsNever.Range('B14').Value = 'I''ve never bought a perian rug for my bathroom'
$wsNever.Range('A14').Value = 'When did you last buy a persian rug for your bathroom?'
$wsNever.Range('B16').Value = 'I''ve never danced flamenco'
$wsWarm.Range('B480').Value = 'In a shocking turn of events'
$wsWarm.Range('A480').Value = 'en un asombroso giro de los acontecimientos'
$wsWarm.Range('C480').Value = 'expressions'
$wsWarm.Range('A481').Value = 'Casi nade sabie mi nombre real'
$wsWarm.Range('B481').Value = 'Hardly anyone knows my real name / almost no one knows my real name'
$wsWarm.Range('A482').Value = 'Casi todo el mundo se quedó dormido'
$wsWarm.Range('B482').Value = 'Almost everyone fell asleep'
$wsWarm.Range('C481').Value = 'casi nadie/casi todo el mundo'
$wsWarm.Range('C482').Value = 'casi nadie/casi todo el mundo'
$wsWarm.Range('A483').Value = 'los dos se quedaron dormidos'
$wsWarm.Range('A484').Value = 'calmar la sed'
$wsWarm.Range('B484').Value = 'to quench my thirst'
$wsWarm.Range('C484').Value = 'all'
$wsNever.Range('A16').Value = 'When was the las time you danced flamenco?'
$wsWarm.Range('B483').Value = 'Both of them fell asleep'
$wsWarm.Range('C483').Value = 'all'

# --- Update the external-data defined name to cover the new rows ---
foreach ($n in $wb.Names) {
    if ($n.Name -eq "warm!DatosExternos_1") {
        $n.RefersTo = "=warm!`$A`$1:`$B`$484"
    }
}

# --- Move the active tab/selection from "never" to "warm" ---
$wsNever.Range('A5').Select() | Out-Null
$wsWarm.Activate() | Out-Null
$wsWarm.Range('B485').Select() | Out-Null
